$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the projection years (columns L:AQ, i.e. 2019-2050) are cleared,
# leaving only the base year (column K, 2018) populated.
$rowsToClear = @(18, 19, 20, 22, 23, 24, 25, 27, 28, 30, 31)
foreach ($r in $rowsToClear) {
    $rng = $ws.Range("L" + $r + ":AQ" + $r)
    $rng.ClearContents()
}

# Rows where the projection years (columns L:AQ) are instead flattened to
# repeat the base year (column K) value across the whole horizon.
$rowsToFlatten = @(26, 29)
foreach ($r in $rowsToFlatten) {
    $baseValue = $ws.Range("K" + $r).Value2
    $rng = $ws.Range("L" + $r + ":AQ" + $r)
    $rng.Value2 = $baseValue
}
